$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("posts")
$ws.Activate()

# Row 11 (post id 10) gets a new generated image + matching color swatch,
# replacing the old "cardtoon..." image / "#70BB9A" color pair.
$ws.Range("E11").Value = "kid-petting-a-dog-with-yellow-background.jpeg"
$ws.Range("F11").Value = "#FCF2D7"

$ws.Range("F13").Select()
